$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I98").Value = 17858412
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 17858412
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -17856914
$ws.Range("N98").ClearContents()
$ws.Range("I122").Value = 17858412
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 53575236
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -53572786
$ws.Range("N122").ClearContents()
$ws.Range("H138").Value = 2435.68
$ws.Range("I138").Value = 1138.0667
$ws.Range("J138").Value = 2664.6707
$ws.Range("K138").Value = 3414.2001
$ws.Range("L138").Value = 7994.0121
$ws.Range("M138").Value = 1725.7999
$ws.Range("N138").Value = -18274.0121
$ws.Range("H141").Value = 2080
$ws.Range("I141").Value = 1446.1364
$ws.Range("J141").Value = 3347.7273
$ws.Range("K141").Value = 4338.4092
$ws.Range("L141").Value = 10043.1819
$ws.Range("M141").Value = 841.5907999999999
$ws.Range("N141").Value = -20403.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 15604.9
$ws.Range("J44").Value = 15604.9
$ws.Range("L44").Value = 15604.9
$ws.Range("N44").Value = -16580.9
$ws.Range("H55").Value = 11000
$ws.Range("J55").Value = 11333.333
$ws.Range("L55").Value = 11333.333
$ws.Range("N55").Value = -11963.333
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 3000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 15000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -21864
$ws.Range("H80").Value = 15179.23
$ws.Range("J80").Value = 15179.23
$ws.Range("L80").Value = 15179.23
$ws.Range("N80").Value = -17175.23
$ws.Range("H83").Value = 15179.23
$ws.Range("J83").Value = 15179.23
$ws.Range("L83").Value = 45537.69
$ws.Range("N83").Value = -55521.69
$ws.Range("H135").Value = 71464.5
$ws.Range("J135").Value = 71464.5
$ws.Range("L135").Value = 71464.5
$ws.Range("N135").Value = -81604.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 49041.773
$ws.Range("I134").Value = 2736.3333
$ws.Range("J134").Value = 148267.72
$ws.Range("K134").Value = 8208.999899999999
$ws.Range("L134").Value = 444803.16
$ws.Range("M134").Value = -5673.999899999999
$ws.Range("N134").Value = -449873.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1037.2
$ws.Range("I31").Value = 871.4
$ws.Range("J31").Value = 1239.8445
$ws.Range("K31").Value = 871.4
$ws.Range("L31").Value = 1239.8445
$ws.Range("M31").Value = -576.4
$ws.Range("N31").Value = -1829.8445
$ws.Range("H34").Value = 1037.2
$ws.Range("I34").Value = 871.4
$ws.Range("J34").Value = 1239.8445
$ws.Range("K34").Value = 871.4
$ws.Range("L34").Value = 1239.8445
$ws.Range("M34").Value = -669.4
$ws.Range("N34").Value = -1643.8445
$ws.Range("H59").Value = 13111.111
$ws.Range("J59").Value = 13111.111
$ws.Range("L59").Value = 13111.111
$ws.Range("N59").Value = -15401.111
$ws.Range("H60").Value = 4321.5
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H68").Value = 29650
$ws.Range("J68").Value = 29650
$ws.Range("L68").Value = 29650
$ws.Range("N68").Value = -31148
$ws.Range("H70").Value = 29560
$ws.Range("J70").Value = 29560
$ws.Range("L70").Value = 29560
$ws.Range("N70").Value = -30190
$ws.Range("H71").Value = 29650
$ws.Range("J71").Value = 29650
$ws.Range("L71").Value = 88950
$ws.Range("N71").Value = -96438
$ws.Range("H73").Value = 29560
$ws.Range("J73").Value = 29560
$ws.Range("L73").Value = 29560
$ws.Range("N73").Value = -31744
$ws.Range("H74").Value = 13111.111
$ws.Range("J74").Value = 13111.111
$ws.Range("L74").Value = 13111.111
$ws.Range("N74").Value = -14859.111
$ws.Range("H77").Value = 13111.111
$ws.Range("J77").Value = 13111.111
$ws.Range("L77").Value = 39333.333
$ws.Range("N77").Value = -48069.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3041.491
$ws.Range("I68").Value = 5902.684
$ws.Range("J68").Value = 1531.4166
$ws.Range("K68").Value = 17708.052
$ws.Range("L68").Value = 4594.2498
$ws.Range("M68").Value = -16897.052
$ws.Range("N68").Value = -6216.2498
$ws.Range("H71").Value = 3041.491
$ws.Range("I71").Value = 5902.684
$ws.Range("J71").Value = 1531.4166
$ws.Range("K71").Value = 53124.156
$ws.Range("L71").Value = 13782.7494
$ws.Range("M71").Value = -49068.156
$ws.Range("N71").Value = -21894.7494
$ws.Range("H113").Value = 1048.6875
$ws.Range("I113").Value = 378.5
$ws.Range("J113").Value = 2165.6667
$ws.Range("K113").Value = 1135.5
$ws.Range("L113").Value = 6497.000100000001
$ws.Range("M113").Value = 1034.5
$ws.Range("N113").Value = -10837.0001
$ws.Range("H129").Value = 26369240
$ws.Range("I129").Value = 476.25
$ws.Range("J129").Value = 45546524
$ws.Range("K129").Value = 1428.75
$ws.Range("L129").Value = 136639572
$ws.Range("M129").Value = 3571.25
$ws.Range("N129").Value = -136649572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 6331.6665
$ws.Range("J57").Value = 6331.6665
$ws.Range("L57").Value = 6331.6665
$ws.Range("N57").Value = -7971.6665
$ws.Range("H62").Value = 47500
$ws.Range("J62").Value = 47500
$ws.Range("L62").Value = 47500
$ws.Range("N62").Value = -48872
$ws.Range("H65").Value = 47500
$ws.Range("J65").Value = 47500
$ws.Range("L65").Value = 142500
$ws.Range("N65").Value = -149364
$ws.Range("H80").Value = 7477.684
$ws.Range("I80").Value = 2406.25
$ws.Range("J80").Value = 11166
$ws.Range("K80").Value = 2406.25
$ws.Range("L80").Value = 11166
$ws.Range("M80").Value = -1408.25
$ws.Range("N80").Value = -13162
$ws.Range("H83").Value = 7477.684
$ws.Range("I83").Value = 2406.25
$ws.Range("J83").Value = 11166
$ws.Range("K83").Value = 12031.25
$ws.Range("L83").Value = 55830
$ws.Range("M83").Value = -7039.25
$ws.Range("N83").Value = -65814

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360

